$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.694.82'
$ws.Range("E2").Value = '  -3.85%  '
$ws.Range("D3").Value = '3.017.85'
$ws.Range("E3").Value = '  -3.20%  '
$ws.Range("D5").Value = "'549.03"
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").Value = "'134.72"
$ws.Range("E6").Value = '  -4.55%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("D8").Value = '3.006.02'
$ws.Range("E8").Value = '  -3.45%  '
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").Value = "'0.148"
$ws.Range("E10").Value = '  -6.32%  '
$ws.Range("D11").Value = "'6.06"
$ws.Range("E11").Value = '  -7.65%  '
$ws.Range("D12").Value = "'0.449"
$ws.Range("E12").Value = '  -3.04%  '
$ws.Range("D13").Value = "'0.0000221"
$ws.Range("E13").Value = '  -3.34%  '
$ws.Range("D14").Value = "'34.29"
$ws.Range("E14").Value = '  -2.78%  '
$ws.Range("D15").Value = '3.510.61'
$ws.Range("E15").Value = '  -3.11%  '
$ws.Range("D16").Value = '61.831.58'
$ws.Range("E16").Value = '  -3.69%  '
$ws.Range("E17").Value = '  -2.87%  '
$ws.Range("D18").Value = '3.022.77'
$ws.Range("E18").Value = '  -2.87%  '
$ws.Range("D19").Value = "'6.65"
$ws.Range("E19").Value = '  -1.67%  '
$ws.Range("D20").Value = "'472.54"
$ws.Range("E20").Value = '  -3.31%  '
$ws.Range("D21").Value = "'13.25"
$ws.Range("E21").Value = '  -2.21%  '
$ws.Range("D22").Value = "'0.674"
$ws.Range("E22").Value = '  -4.97%  '
$ws.Range("D23").Value = "'7.08"
$ws.Range("E23").Value = '  -2.08%  '
$ws.Range("D24").Value = "'80.17"
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("D25").Value = "'12.09"
$ws.Range("E25").Value = '  -2.83%  '
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("D27").Value = "'2.72"
$ws.Range("E27").Value = '  -0.90%  '
$ws.Range("D28").Value = "'7.79"
$ws.Range("E28").Value = '  -5.81%  '
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("E30").Value = '  -1.32%  '
$ws.Range("B31").Value = 'Mantle'
$ws.Range("C31").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'25.71"
$ws.Range("E32").Value = '  -3.31%  '
$ws.Range("D33").Value = "'2.30"
$ws.Range("E33").Value = '  -4.16%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = "'5.46"
$ws.Range("E34").Value = '  +0.76%  '
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").Value = "'55.31"
$ws.Range("E35").Value = '  -3.80%  '
$ws.Range("D36").Value = "'5.90"
$ws.Range("E36").Value = '  -3.06%  '
$ws.Range("D37").Value = "'458.26"
$ws.Range("E37").Value = '  -9.58%  '
$ws.Range("D38").Value = '3.217.63'
$ws.Range("E38").Value = '  -1.84%  '
$ws.Range("D39").Value = "'0.0797"
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("D40").Value = "'0.0383"
$ws.Range("E40").Value = '  -6.05%  '
$ws.Range("D41").Value = "'0.118"
$ws.Range("E41").Value = '  -1.54%  '
$ws.Range("D42").Value = "'8.15"
$ws.Range("E42").Value = '  -0.55%  '
$ws.Range("E43").Value = '  -11.14%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").Value = "'26.02"
$ws.Range("E45").Value = '  +3.36%  '
$ws.Range("D46").Value = "'0.244"
$ws.Range("E46").Value = '  -5.88%  '
$ws.Range("D47").Value = "'1.99"
$ws.Range("E47").Value = '  -4.59%  '
$ws.Range("E48").Value = '  -1.22%  '
$ws.Range("D49").Value = "'118.03"
$ws.Range("E49").Value = '  -4.40%  '
$ws.Range("D50").Value = '0.0₃0495'
$ws.Range("E50").Value = '  -8.56%  '
